# Apply 2025-07-21 violent crime data updates across affected worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3689
$ws.Range("L3").Value = 3844
$ws.Range("L4").Value = 954
$ws.Range("L5").Value = 231
$ws.Range("L6").Value = 3364
$ws.Range("L7").Value = 12082

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 228
$ws.Range("L3").Value = 260
$ws.Range("L6").Value = 214
$ws.Range("L7").Value = 781

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 183
$ws.Range("L7").Value = 564

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 171

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 133
$ws.Range("L6").Value = 122
$ws.Range("L7").Value = 436

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 84
$ws.Range("L4").Value = 17
$ws.Range("L7").Value = 202

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 100
$ws.Range("L3").Value = 10
$ws.Range("L6").Value = 96
$ws.Range("L8").Value = 781
$ws.Range("L11").Value = 203
$ws.Range("L19").Value = 341
$ws.Range("L20").Value = 304
$ws.Range("L26").Value = 14
$ws.Range("L27").Value = 110
$ws.Range("L29").Value = 660
$ws.Range("L31").Value = 118
$ws.Range("L33").Value = 564
$ws.Range("L34").Value = 74
$ws.Range("L35").Value = 19
$ws.Range("L36").Value = 162
$ws.Range("L37").Value = 436
$ws.Range("L42").Value = 384
$ws.Range("L43").Value = 90
$ws.Range("L47").Value = 88
$ws.Range("L51").Value = 149
$ws.Range("L52").Value = 246
$ws.Range("L53").Value = 136
$ws.Range("L54").Value = 251
$ws.Range("L59").Value = 17
$ws.Range("L60").Value = 73
$ws.Range("L63").Value = 38
$ws.Range("L64").Value = 79
$ws.Range("L67").Value = 430
$ws.Range("L68").Value = 39
$ws.Range("L69").Value = 31
$ws.Range("L71").Value = 32
$ws.Range("L76").Value = 180
$ws.Range("L78").Value = 152
$ws.Range("L79").Value = 316
$ws.Range("L85").Value = 630
$ws.Range("L88").Value = 135
$ws.Range("L89").Value = 170
$ws.Range("L90").Value = 114
$ws.Range("L91").Value = 175
$ws.Range("L94").Value = 151
$ws.Range("L95").Value = 171
$ws.Range("L99").Value = 202
$ws.Range("L101").Value = 12082

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 42
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 123
$ws.Range("L3").Value = 165
$ws.Range("L7").Value = 430

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 59
$ws.Range("L6").Value = 118
$ws.Range("L7").Value = 251

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 201
$ws.Range("L3").Value = 249
$ws.Range("L6").Value = 167
$ws.Range("L7").Value = 660

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 101
$ws.Range("L7").Value = 341

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L4").Value = 23
$ws.Range("L6").Value = 82
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 41
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 96

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 110
$ws.Range("L7").Value = 384

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 152

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 62
$ws.Range("L3").Value = 74
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 175

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 110
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 316

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 96
$ws.Range("L7").Value = 304

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 62
$ws.Range("L7").Value = 162

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 36
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 151

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 31
$ws.Range("L3").Value = 31
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("L5").Value = 9
$ws.Range("L6").Value = 14

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 77
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 49
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 170

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 149

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 184
$ws.Range("L3").Value = 258
$ws.Range("L6").Value = 132
$ws.Range("L7").Value = 630

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 10

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 84
$ws.Range("L3").Value = 76
$ws.Range("L6").Value = 64
$ws.Range("L7").Value = 246
